# Bug fix: the "Level" column (B) was left blank for rows 59-100 because the
# ratio formula (Price / Closing1d) was never entered for those rows, and the
# column's number format was wrongly set to a percentage instead of a plain
# decimal number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fill in the missing "Level" ratio formula for rows 59 through 100.
for ($row = 59; $row -le 100; $row++) {
    $ws.Range("B$row").Formula = "=F$row/G$row"
}

# 2) Fix the number format applied to those same cells: it should show a
#    plain decimal (0.000) rather than a percentage (0.000%).
$ws.Range("B59:B100").NumberFormat = "0.000"
